# SGU Science or Fiction - add episode 602 results
$wb = $excel.ActiveWorkbook

$wsData    = $wb.Worksheets.Item("Data")
$wsResults = $wb.Worksheets.Item("Results")
$wsSummary = $wb.Worksheets.Item("Summary")

# ---------------------------------------------------------------
# 1. Data sheet - new row for episode 602
# ---------------------------------------------------------------
$wsData.Range("A4").Value = 602
$wsData.Range("B4").Value = "Global warming"
$wsData.Range("C4").Value = "Carl Sagan was the first scientist to publicly warn about the possibility of manmade global warming from greenhouse gas emissions, in a 1980 essay"
$wsData.Range("D4").Value = "The 15 hottest years on record since 1880 have all been since 1998. "
$wsData.Range("E4").Value = "Climate models show that even if CO2 emissions were stopped entirely, global temperatures would continue to rise for at least a century"
$wsData.Range("F4").Formula = "=NA()"
$wsData.Range("G4").Value = 1
$wsData.Range("H4").Value = "Steve"
$wsData.Range("I4").Value = 2
$wsData.Range("J4").Value = 1
$wsData.Range("K4").Value = 1
$wsData.Range("L4").Value = 1
$wsData.Range("M4").Formula = "=NA()"
$wsData.Range("N4").Formula = "=NA()"

# Centered style for the new data row plus pre-formatted, still-empty
# tracker cells for upcoming episodes
$wsData.Range("G4:L53").HorizontalAlignment = -4108

# ---------------------------------------------------------------
# 2. Results sheet - mirror row for episode 602
# ---------------------------------------------------------------
$wsResults.Range("A4").Formula = "=Data!A4"
$wsResults.Range("B4").Formula = "=Data!B4"
$wsResults.Range("C4").Formula = "=Data!H4"
$wsResults.Range("D4").Formula = "=IF(Data!I4=Data!`$G4,1,0)"
$wsResults.Range("E4").Formula = "=IF(Data!J4=Data!`$G4,1,0)"
$wsResults.Range("F4").Formula = "=IF(Data!K4=Data!`$G4,1,0)"
$wsResults.Range("G4").Formula = "=IF(Data!L4=Data!`$G4,1,0)"
$wsResults.Range("H4").Formula = "=IF(Data!M4=Data!`$G4,1,0)"
$wsResults.Range("I4").Formula = "=IF(Data!N4=Data!`$G4,1,0)"
$wsResults.Range("J4").Formula = '=IF(SUMIF(D4:I4,"<>#N/A")=0, 1, 0)'

# ---------------------------------------------------------------
# 3. Summary sheet - add leading "Overall % Wins" label column
# ---------------------------------------------------------------
$wsSummary.Columns.Item(1).Insert()
$wsSummary.Columns.Item(1).ColumnWidth = 25.1666667
$wsSummary.Range("A2").Value = "Overall % Wins"

# ---------------------------------------------------------------
# 4. Selections / active sheet
# ---------------------------------------------------------------
$wsSummary.Range("B2").Select()
$wsResults.Range("A4:J4").Select()
$wsData.Range("A5").Select()
$wsData.Activate()
